$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells.Item(28, 8).Value = 2108.25
$ws.Cells.Item(28, 9).Value = 1876.4
$ws.Cells.Item(28, 11).Value = 1876.4
$ws.Cells.Item(28, 13).Value = -1391.4
# Row 38
$ws.Cells.Item(38, 8).Value = 21.166666
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 14).ClearContents()
# Row 58
$ws.Cells.Item(58, 8).Value = 2078.75
$ws.Cells.Item(58, 10).Value = 6200
$ws.Cells.Item(58, 12).Value = 18600
$ws.Cells.Item(58, 14).Value = -18900
# Row 97
$ws.Cells.Item(97, 8).Value = 875
$ws.Cells.Item(97, 10).Value = 875
$ws.Cells.Item(97, 12).Value = 2625
$ws.Cells.Item(97, 14).Value = -3617
# Row 107
$ws.Cells.Item(107, 8).Value = 1238.0435
$ws.Cells.Item(107, 9).Value = 904.35297
$ws.Cells.Item(107, 11).Value = 904.35297
$ws.Cells.Item(107, 13).Value = 1015.64703
# Row 113
$ws.Cells.Item(113, 8).Value = 6416.7393
$ws.Cells.Item(113, 10).Value = 6366.8
$ws.Cells.Item(113, 12).Value = 6366.8
$ws.Cells.Item(113, 14).Value = -12874.8
# Row 132
$ws.Cells.Item(132, 8).Value = 3191555.2
$ws.Cells.Item(132, 9).Value = 3419438
$ws.Cells.Item(132, 11).Value = 10258314
$ws.Cells.Item(132, 13).Value = -10255784

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 3852.923
$ws.Cells.Item(45, 9).Value = 2815.3333
$ws.Cells.Item(45, 10).Value = 4164.2
$ws.Cells.Item(45, 11).Value = 2815.3333
$ws.Cells.Item(45, 12).Value = 4164.2
$ws.Cells.Item(45, 13).Value = -2438.3333
$ws.Cells.Item(45, 14).Value = -4918.2
# Row 61
$ws.Cells.Item(61, 8).Value = 10994
$ws.Cells.Item(61, 9).Value = 1104.4445
$ws.Cells.Item(61, 10).Value = 100000
$ws.Cells.Item(61, 11).Value = 1104.4445
$ws.Cells.Item(61, 12).Value = 100000
$ws.Cells.Item(61, 13).Value = -892.4445000000001
$ws.Cells.Item(61, 14).Value = -100424
# Row 74
$ws.Cells.Item(74, 8).Value = 612392.4
$ws.Cells.Item(74, 9).Value = 1500877.5
$ws.Cells.Item(74, 10).Value = 20069
$ws.Cells.Item(74, 11).Value = 1500877.5
$ws.Cells.Item(74, 12).Value = 20069
$ws.Cells.Item(74, 13).Value = -1500003.5
$ws.Cells.Item(74, 14).Value = -21817
# Row 77
$ws.Cells.Item(77, 8).Value = 612392.4
$ws.Cells.Item(77, 9).Value = 1500877.5
$ws.Cells.Item(77, 10).Value = 20069
$ws.Cells.Item(77, 11).Value = 7504387.5
$ws.Cells.Item(77, 12).Value = 100345
$ws.Cells.Item(77, 13).Value = -7500019.5
$ws.Cells.Item(77, 14).Value = -109081
# Row 97
$ws.Cells.Item(97, 8).Value = 746
$ws.Cells.Item(97, 9).Value = 663.2093
$ws.Cells.Item(97, 10).Value = 1932.6666
$ws.Cells.Item(97, 11).Value = 663.2093
$ws.Cells.Item(97, 12).Value = 1932.6666
$ws.Cells.Item(97, 13).Value = -167.2093
$ws.Cells.Item(97, 14).Value = -2924.6666
# Row 132
$ws.Cells.Item(132, 8).Value = 1419.1072
$ws.Cells.Item(132, 10).Value = 3495.7144
$ws.Cells.Item(132, 12).Value = 10487.1432
$ws.Cells.Item(132, 14).Value = -15547.1432
# Row 136
$ws.Cells.Item(136, 8).Value = 10994
$ws.Cells.Item(136, 9).Value = 1104.4445
$ws.Cells.Item(136, 10).Value = 100000
$ws.Cells.Item(136, 11).Value = 3313.3335
$ws.Cells.Item(136, 12).Value = 300000
$ws.Cells.Item(136, 13).Value = -763.3335000000002
$ws.Cells.Item(136, 14).Value = -305100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 1210.44
$ws.Cells.Item(20, 9).Value = 965.7857
$ws.Cells.Item(20, 10).Value = 1521.8182
$ws.Cells.Item(20, 11).Value = 965.7857
$ws.Cells.Item(20, 12).Value = 1521.8182
$ws.Cells.Item(20, 13).Value = -718.7857
$ws.Cells.Item(20, 14).Value = -2015.8182
# Row 94
$ws.Cells.Item(94, 8).Value = 5921.778
$ws.Cells.Item(94, 9).Value = 5921.778
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 5921.778
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = -5470.778
$ws.Cells.Item(94, 14).ClearContents()
# Row 99
$ws.Cells.Item(99, 8).Value = 1500
$ws.Cells.Item(99, 9).Value = 1500
$ws.Cells.Item(99, 11).Value = 1500
$ws.Cells.Item(99, 13).Value = -2
# Row 134
$ws.Cells.Item(134, 8).Value = 1545.5
$ws.Cells.Item(134, 9).Value = 1545.5
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 4636.5
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = -2101.5
$ws.Cells.Item(134, 14).ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 1900.5555
$ws.Cells.Item(16, 9).Value = 1814.0667
$ws.Cells.Item(16, 11).Value = 1814.0667
$ws.Cells.Item(16, 13).Value = -1527.0667
# Row 26
$ws.Cells.Item(26, 8).Value = 30277.5
$ws.Cells.Item(26, 9).Value = 5000
$ws.Cells.Item(26, 10).Value = 55555
$ws.Cells.Item(26, 11).Value = 5000
$ws.Cells.Item(26, 12).Value = 55555
$ws.Cells.Item(26, 13).Value = -4713
$ws.Cells.Item(26, 14).Value = -56129
# Row 58
$ws.Cells.Item(58, 8).Value = 2052.5
$ws.Cells.Item(58, 9).Value = 1161
$ws.Cells.Item(58, 10).Value = 4727
$ws.Cells.Item(58, 11).Value = 1161
$ws.Cells.Item(58, 12).Value = 4727
$ws.Cells.Item(58, 13).Value = -958
$ws.Cells.Item(58, 14).Value = -5133
# Row 93
$ws.Cells.Item(93, 8).Value = 29333.334
$ws.Cells.Item(93, 9).Value = 29333.334
$ws.Cells.Item(93, 11).Value = 29333.334
$ws.Cells.Item(93, 13).Value = -27461.334
# Row 94
$ws.Cells.Item(94, 8).Value = 1568.5238
$ws.Cells.Item(94, 10).Value = 1922.6
$ws.Cells.Item(94, 12).Value = 1922.6
$ws.Cells.Item(94, 14).Value = -2824.6
# Row 105
$ws.Cells.Item(105, 8).Value = 1090.875
$ws.Cells.Item(105, 9).Value = 830.26666
$ws.Cells.Item(105, 10).Value = 5000
$ws.Cells.Item(105, 11).Value = 830.26666
$ws.Cells.Item(105, 12).Value = 5000
$ws.Cells.Item(105, 13).Value = 916.73334
$ws.Cells.Item(105, 14).Value = -8494
# Row 113
$ws.Cells.Item(113, 8).Value = 1900.5555
$ws.Cells.Item(113, 9).Value = 1814.0667
$ws.Cells.Item(113, 11).Value = 1814.0667
$ws.Cells.Item(113, 13).Value = 355.9332999999999
# Row 136
$ws.Cells.Item(136, 8).Value = 2052.5
$ws.Cells.Item(136, 9).Value = 1161
$ws.Cells.Item(136, 10).Value = 4727
$ws.Cells.Item(136, 11).Value = 3483
$ws.Cells.Item(136, 12).Value = 14181
$ws.Cells.Item(136, 13).Value = -933
$ws.Cells.Item(136, 14).Value = -19281

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Cells.Item(3, 8).Value = 959.5714
$ws.Cells.Item(3, 9).Value = 959.5714
$ws.Cells.Item(3, 11).Value = 2878.7142
$ws.Cells.Item(3, 13).Value = -2766.7142
# Row 24
$ws.Cells.Item(24, 8).Value = 1999
$ws.Cells.Item(24, 10).Value = 1999
$ws.Cells.Item(24, 12).Value = 5997
$ws.Cells.Item(24, 14).Value = -6457
# Row 64
$ws.Cells.Item(64, 8).Value = 7500
$ws.Cells.Item(64, 9).Value = 50000
$ws.Cells.Item(64, 11).Value = 150000
$ws.Cells.Item(64, 13).Value = -149730
# Row 67
$ws.Cells.Item(67, 8).Value = 7500
$ws.Cells.Item(67, 9).Value = 50000
$ws.Cells.Item(67, 11).Value = 150000
$ws.Cells.Item(67, 13).Value = -149064
# Row 69
$ws.Cells.Item(69, 8).Value = 4466.875
$ws.Cells.Item(69, 10).Value = 4964.2856
$ws.Cells.Item(69, 12).Value = 14892.8568
$ws.Cells.Item(69, 14).Value = -16514.8568
# Row 72
$ws.Cells.Item(72, 8).Value = 4466.875
$ws.Cells.Item(72, 10).Value = 4964.2856
$ws.Cells.Item(72, 12).Value = 44678.5704
$ws.Cells.Item(72, 14).Value = -52790.5704
# Row 103
$ws.Cells.Item(103, 8).Value = 1131.25
$ws.Cells.Item(103, 9).Value = 1131.25
$ws.Cells.Item(103, 11).Value = 3393.75
$ws.Cells.Item(103, 13).Value = -2514.75
# Row 107
$ws.Cells.Item(107, 8).Value = 1320
$ws.Cells.Item(107, 9).Value = 1615.6428
$ws.Cells.Item(107, 11).Value = 4846.928400000001
$ws.Cells.Item(107, 13).Value = -2926.928400000001
# Row 124
$ws.Cells.Item(124, 8).Value = 9074.727999999999
$ws.Cells.Item(124, 10).Value = 9364
$ws.Cells.Item(124, 12).Value = 28092
$ws.Cells.Item(124, 14).Value = -37912
# Row 129
$ws.Cells.Item(129, 8).Value = 2501.7646
$ws.Cells.Item(129, 10).Value = 2616.875
$ws.Cells.Item(129, 12).Value = 7850.625
$ws.Cells.Item(129, 14).Value = -17850.625
# Row 131
$ws.Cells.Item(131, 8).Value = 226506
$ws.Cells.Item(131, 10).Value = 1721.4
$ws.Cells.Item(131, 12).Value = 5164.200000000001
$ws.Cells.Item(131, 14).Value = -15244.2
# Row 134
$ws.Cells.Item(134, 8).Value = 450
$ws.Cells.Item(134, 9).Value = 450
$ws.Cells.Item(134, 11).Value = 1350
$ws.Cells.Item(134, 13).Value = 3720

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 18533.275
$ws.Cells.Item(102, 9).Value = 21008.6
$ws.Cells.Item(102, 11).Value = 21008.6
$ws.Cells.Item(102, 13).Value = -19386.6
# Row 132
$ws.Cells.Item(132, 8).Value = 1622.5
$ws.Cells.Item(132, 10).Value = 2707
$ws.Cells.Item(132, 12).Value = 8121
$ws.Cells.Item(132, 14).Value = -13181

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 111
$ws.Cells.Item(111, 8).Value = 59999
$ws.Cells.Item(111, 10).Value = 59999
$ws.Cells.Item(111, 12).Value = 59999
$ws.Cells.Item(111, 14).Value = -68179
# Row 132
$ws.Cells.Item(132, 8).Value = 10428.286
$ws.Cells.Item(132, 9).Value = 12199.6
$ws.Cells.Item(132, 11).Value = 36598.8
$ws.Cells.Item(132, 13).Value = -34068.8
# Row 136
$ws.Cells.Item(136, 8).Value = 3000.25
$ws.Cells.Item(136, 9).Value = 2914.2354
$ws.Cells.Item(136, 11).Value = 8742.706200000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Cells.Item(122, 8).Value = 23946.746
$ws.Cells.Item(122, 9).Value = 26058
$ws.Cells.Item(122, 11).Value = 78174
$ws.Cells.Item(122, 13).Value = -75724
# Row 126
$ws.Cells.Item(126, 8).Value = 220746.05
$ws.Cells.Item(126, 9).Value = 2902.8948
$ws.Cells.Item(126, 11).Value = 8708.6844
$ws.Cells.Item(126, 13).Value = -6238.6844
